$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# Fix / de-duplicate a few strings in the "University" column (trailing
# whitespace + capitalization/typo fixes) so they collapse onto existing
# shared-string entries instead of keeping their own duplicate entries.
$ws.Range("C4").Value = "Southern Methodist University"
$ws.Range("C6").Value = "University of North Texas"
$ws.Range("C16").Value = "Centennial High School"
$ws.Range("C17").Value = "Centennial High School"

# Update the active selection to reflect where the author was last working.
$null = $ws.Range("C24").Select()
